# "se actualiza el diccionario de datos" (the data dictionary is updated)
#
# The source workbook is a one-row "data dictionary" extract: row 1 holds the
# field names, row 2 holds a sample record. Column J's field was renamed from
# "Nit_tercero" to "identificacion_tercero" (the ID-of-the-third-party field
# in the dictionary got relabeled). Everything else in that column/row is
# untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "identificacion_tercero"

# Leave the selection where the author apparently ended up after the edit.
$ws.Range("P1").Select()
